# Changes of 15th july 2022
# Update the FedEx shipment tracking numbers in column P (rows 2-26) with
# the next batch of tracking numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$trackingNumbers = @(
  "320017965214",
  "320017965225",
  "320017965258",
  "320017965270",
  "320017965317",
  "320017965339",
  "320017965361",
  "320017965383",
  "320017965410",
  "320017965431",
  "320017965475",
  "320017965497",
  "320017965523",
  "320017965545",
  "320017965578",
  "320017965590",
  "320017965637",
  "320017965659",
  "320017965681",
  "320017965707",
  "320017965730",
  "320017965740",
  "320017965751",
  "320017965762",
  "320017965773"
)

for ($i = 0; $i -lt $trackingNumbers.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 16)
    # Leading apostrophe forces Excel to store the value as text (keeping
    # the long digit string intact instead of coercing it to a number).
    $cell.Value = "'" + $trackingNumbers[$i]
    # Reset to the default "Normal" style so we don't leave the cell with
    # a quote-prefix / text number-format style that the original file
    # never had (column P cells here carry no explicit style).
    $cell.Style = "Normal"
}
